# TASK 01.xlsx — "Add files via upload / Added Task 1: Gender distribution
# chart in Excel" follow-up edit:
#   - rename the third GENDER category from "       OTHER" to "NON-BINARY"
#   - leave the selection on B8 (next empty row under the data) instead of Q12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The gender labels live in B5:B7 (MALE, FEMALE, OTHER). The third one is
# being renamed to a plain "NON-BINARY" (no leading spaces / preserved
# whitespace like the other labels use).
$ws.Range("B7").Value = "NON-BINARY"

# Move/leave the active selection at B8, just under the table.
$ws.Range("B8").Select() | Out-Null
